$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the styled/bordered "data row" formatting (currently s="2", same as
# row 16) down through row 21 so the new rows pick up the existing style
# instead of minting a new cellXf.
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C21").PasteSpecial(-4122)

# Row 11: ClaimMojio / Done / N  (was ClaimMojio / <blank> / Y)
$ws.Range("A11").Value = "ClaimMojio"
$ws.Range("B11").Value = "Done"
$ws.Range("C11").Value = "N"

# Row 12: BillingInformation / Client stopped it / N
$ws.Range("A12").Value = "BillingInformation"
$ws.Range("B12").Value = "Client stopped it"
$ws.Range("C12").Value = "N"

# Row 13: Contact_Support / Done / N
$ws.Range("A13").Value = "Contact_Support"
$ws.Range("B13").Value = "Done"
$ws.Range("C13").Value = "N"

# Row 14: MojioShop / Issue in Details link in firefox / N
$ws.Range("A14").Value = "MojioShop"
$ws.Range("B14").Value = "Issue in Details link in firefox"
$ws.Range("C14").Value = "N"

# Row 15: MyOrders / <blank> / Y
$ws.Range("A15").Value = "MyOrders"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "Y"

# Row 16: LocateMojio / Done / N
$ws.Range("A16").Value = "LocateMojio"
$ws.Range("B16").Value = "Done"
$ws.Range("C16").Value = "N"

# Row 17: TripHistory / Date search remains / N
$ws.Range("A17").Value = "TripHistory"
$ws.Range("B17").Value = "Date search remains"
$ws.Range("C17").Value = "N"

# Row 18: SendFeedback / Done / N
$ws.Range("A18").Value = "SendFeedback"
$ws.Range("B18").Value = "Done"
$ws.Range("C18").Value = "N"

# Row 19: Settings / <blank> / N
$ws.Range("A19").Value = "Settings"
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = "N"

# Rows 20-21: fully blank, styled rows (extends the used range to A1:C21).
$ws.Range("A20").Value = ""
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("A21").Value = ""
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = ""

# Move the active selection to C15, matching the saved view state.
$ws.Range("C15").Select()
